# Mise à jour des questions
# Adds 18 new trivia question rows (rows 43-60) to the "Feuil1" worksheet,
# matching the quiz format already used in the sheet (Question, 4 possible
# answers, and the correct answer repeated in column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = "Combien d'albums un chanteur doit-il vendre pour recevoir un disque d'or ?"
$ws.Range("B43").Value = "100'000"
$ws.Range("C43").Value = "70'000"
$ws.Range("D43").Value = "80'000"
$ws.Range("E43").Value = "75'000"
$ws.Range("F43").Value = "75'000"

$ws.Range("A44").Value = "De quel moyen de locomotion le grand bi est-il l'ancêtre ?"
$ws.Range("B44").Value = "Le dirigeable"
$ws.Range("C44").Value = "Le monocycle"
$ws.Range("D44").Value = "Le trycycle"
$ws.Range("E44").Value = "La bicyclette"
$ws.Range("F44").Value = "La bicyclette"

$ws.Range("A45").Value = "Quel film d'animation a pour héros Woody le cow-boy et Buzz l'éclair ?"
$ws.Range("C45").Value = "Cars"
$ws.Range("D45").Value = "Robots"
$ws.Range("B45").Value = "Small Soldiers"
$ws.Range("E45").Value = "Toy Story"
$ws.Range("F45").Value = "Toy Story"

$ws.Range("A46").Value = "Quelle est la particularité du tonneau des Danaïdes ?"
$ws.Range("B46").Value = "Il ne se dilate jamais"
$ws.Range("C46").Value = "Il est sans fond "
$ws.Range("D46").Value = "Il est tout le temps ouvert"
$ws.Range("E46").Value = "Il est en cristal"
$ws.Range("F46").Value = "Il est sans fond"

$ws.Range("A47").Value = "Quel cavalier est le maître de Rossinante ? "
$ws.Range("B47").Value = "Don Quichotte"
$ws.Range("F47").Value = "Don Quichotte"
$ws.Range("C47").Value = "Jorge Luis Borges"
$ws.Range("D47").Value = "Günter Grass"
$ws.Range("E47").Value = "Sirano de Bergerac"

$ws.Range("A48").Value = "Quelle est la capitale du Cameroun ? "
$ws.Range("B48").Value = "Douala"
$ws.Range("C48").Value = "Bertoua"
$ws.Range("D48").Value = "Yaoundé"
$ws.Range("F48").Value = "Yaoundé"
$ws.Range("E48").Value = "Luanda"

$ws.Range("A49").Value = "Qui est la muse de l'Histoire dans la mythologie grecque ?"
$ws.Range("B49").Value = "Clio"
$ws.Range("F49").Value = "Clio"
$ws.Range("C49").Value = "Calliope"
$ws.Range("D49").Value = "Uranie"
$ws.Range("E49").Value = "Euterpe"

$ws.Range("A50").Value = "Dans la mythologie grecque, qui est le dieu des voleurs ? "
$ws.Range("B50").Value = "Horus"
$ws.Range("C50").Value = "Hermès"
$ws.Range("F50").Value = "Hermès"
$ws.Range("D50").Value = "Helheim"
$ws.Range("E50").Value = "Hélios"

$ws.Range("A51").Value = "Quel est le véritable prénom de M. Pokora ? "
$ws.Range("B51").Value = "Matt"
$ws.Range("D51").Value = "Mathéo"
$ws.Range("E51").Value = "Mathias"
$ws.Range("C51").Value = "Matthieu"
$ws.Range("F51").Value = "Matthieu"

$ws.Range("A52").Value = "Quel couple de petits personnages rouge et bleu s'anime en 1974 ? "
$ws.Range("B52").Value = "Boule & Bill"
$ws.Range("C52").Value = "Titi & Grosminet"
$ws.Range("D52").Value = "Coyotte & Bip bip"
$ws.Range("E52").Value = "Chapi & Chapo"
$ws.Range("F52").Value = "Chapi & Chapo"

$ws.Range("A53").Value = "Quel sirop consommez-vous si vous buvez un Monaco ?"
$ws.Range("B53").Value = "Cerise"
$ws.Range("C53").Value = "Framboise"
$ws.Range("D53").Value = "Menthe "
$ws.Range("E53").Value = "Grenadine"
$ws.Range("F53").Value = "Grenadine"

$ws.Range("A54").Value = "Comment s'appelle le chevalier dans Zelda"
$ws.Range("B54").Value = "Link"
$ws.Range("F54").Value = "Link"
$ws.Range("C54").Value = "Zelda"
$ws.Range("D54").Value = "Gooruk"
$ws.Range("E54").Value = "Revali"

$ws.Range("A55").Value = "Comment s'appelle la princesse dans Mario Galaxy"
$ws.Range("C55").Value = "Daisy"
$ws.Range("B55").Value = "Peach"
$ws.Range("D55").Value = "Harmonie"
$ws.Range("F55").Value = "Harmonie"
$ws.Range("E55").Value = "Queen"

$ws.Range("A56").Value = "Dans quoi se cache le célèbre personnage dans Metal gear solid"
$ws.Range("B56").Value = "Un carton"
$ws.Range("F56").Value = "Un carton"
$ws.Range("C56").Value = "Un buisson"
$ws.Range("D56").Value = "Une cabine télephonique"
$ws.Range("E56").Value = "Un sac de couchage"

$ws.Range("A57").Value = "Combien y'a-t-il de pokémons dans la 1er géneration"
$ws.Range("B57").Value = 160
$ws.Range("C57").Value = 151
$ws.Range("D57").Value = 161
$ws.Range("E57").Value = 150
$ws.Range("F57").Value = 151

$ws.Range("A58").Value = "Quel est la mascotte dans pokémon"
$ws.Range("B58").Value = "Pikachu"
$ws.Range("F58").Value = "Pikachu"
$ws.Range("C58").Value = "Evoli"
$ws.Range("D58").Value = "Rondoudou"
$ws.Range("E58").Value = "Insolourdo"

$ws.Range("A59").Value = "Un berger a 16 brebis, toutes meurent sauf 10. Combien lui en reste-t-il ?"
$ws.Range("B59").Value = 10
$ws.Range("C59").Value = 6
$ws.Range("D59").Value = 16
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 10

$ws.Range("A60").Value = "Au Japon, qu'est-ce qu'un yakuza ? "
$ws.Range("B60").Value = "Un vendeur de drogue"
$ws.Range("C60").Value = "Une vendeur de tapis"
$ws.Range("D60").Value = "Un membre de la mafia"
$ws.Range("F60").Value = "Un membre de la mafia"
$ws.Range("E60").Value = "Un pilier économique"

# Move the selection to match the end of the newly-added data (mirrors the
# author's saved window state after entering the new questions).
$ws.Range("F60").Select()
